$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for data rows 2-244 is updated from serial date
# 45179 (2023-09-10) to 45180 (2023-09-11).
$ws.Range("C2:C244").Value = 45180
